$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / "want to go" count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 583
$ws1.Range("F5").Value = 530
$ws1.Range("F6").Value = 298
$ws1.Range("F7").Value = 2686
$ws1.Range("F8").Value = 458
$ws1.Range("F9").Value = 7453
$ws1.Range("F12").Value = 26
$ws1.Range("F13").Value = 234
$ws1.Range("F14").Value = 40

# Sheet "全部类型" (All types) - same underlying data, update corresponding rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 583
$ws4.Range("F5").Value = 530
$ws4.Range("F6").Value = 298
$ws4.Range("F9").Value = 2686
$ws4.Range("F10").Value = 458
$ws4.Range("F11").Value = 7453
$ws4.Range("F14").Value = 26
$ws4.Range("F17").Value = 234
$ws4.Range("F18").Value = 40
